# Generate Report for handoff
#
# Adds two newly-discovered source files
# (2513cad1-8dc2-49f5-aa66-cdb79a8a50da.md and
#  9834a0f2-8b5f-4a90-9762-e6e5262ae73d.md) to the localization-status
# report, pushes the two already-tracked files
# (28cc7299-... and fdc78f58-...) from "Ready for handoff" into
# "In Translation" with refreshed handoff timestamps, and keeps the
# ".localization-config" / "Not to be localized" row trailing the table
# on every sheet.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # BGR packing of RGB(100,149,237) == style "FF6495ED"

function Set-PlainCell {
    param($ws, $cellRef, $text)
    $ws.Range($cellRef).Value = $text
}

function Set-LinkCell {
    param($ws, $cellRef, $text, $url)
    $ws.Range($cellRef).Value = $text
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, $null, $null, $text)
    $ws.Range($cellRef).Font.Underline = $true
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
}

function Set-DateCell {
    param($ws, $cellRef, $text)
    $ws.Range($cellRef).Value = $text
    $ws.Range($cellRef).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# File / hash identifiers
# ---------------------------------------------------------------------
$file28 = "28cc7299-bfbe-42ea-8d00-d1cda1a91de3.md"
$file28hash = "28cc7299-bfbe-42ea-8d00-d1cda1a91de3.f7bf3c05e96e9cd68c6aa5d1351ea639bdad9762"
$filefdc = "fdc78f58-38f3-470a-afd8-03123c183522.md"
$filefdchash = "fdc78f58-38f3-470a-afd8-03123c183522.d736a06a8b9b1cdcd3547941eeb49a0ee5b8b635"
$file2513 = "2513cad1-8dc2-49f5-aa66-cdb79a8a50da.md"
$file2513hash = "2513cad1-8dc2-49f5-aa66-cdb79a8a50da.7ac5dc1129d273edfb69c021d7a686667378652f"
$file9834 = "9834a0f2-8b5f-4a90-9762-e6e5262ae73d.md"
$file9834hash = "9834a0f2-8b5f-4a90-9762-e6e5262ae73d.a2656374f4f56fd5b514f9fc1f5bed2f49b570eb"
$cfgFile = ".localization-config"

$newHandoffZh = "2016-01-13 01:39:03"
$newHandoffDe = "2016-01-13 01:39:25"
$epoch = "0001-01-01 00:00:00"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/7aedd09a01fbf8331e34e14b744b91e7ab361690/e2e"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/7aedd09a01fbf8331e34e14b744b91e7ab361690/.localization-config"
$zhXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/92a1a5d4d1050b3613cf5e56fba6a9b4ae242d0d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang"
$deXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e36ed08ae1e8558c83d72041ce96ebd71abb4405/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang"

# ---------------------------------------------------------------------
# Sheet "Overview" — simple 3 column (File Name / zh-cn / de-de) table
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-LinkCell $wsOverview "A5" $file9834 "$mdBase/$file9834"
Set-PlainCell $wsOverview "B5" "Ready for handoff"
Set-PlainCell $wsOverview "C5" "Ready for handoff"

Set-LinkCell $wsOverview "A6" $cfgFile $cfgUrl
Set-PlainCell $wsOverview "B6" "Not to be localized"
Set-PlainCell $wsOverview "C6" "Not to be localized"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2 / Row 3 — already-tracked files move to "In Translation" with a
# refreshed handoff datetime.
Set-PlainCell $wsZh "B2" "In Translation"
Set-DateCell  $wsZh "D2" $newHandoffZh

Set-PlainCell $wsZh "B3" "In Translation"
Set-DateCell  $wsZh "D3" $newHandoffZh

# Row 4 — new file 2513cad1-...
Set-LinkCell $wsZh "A4" $file2513 "$mdBase/$file2513"
Set-PlainCell $wsZh "B4" "Ready for handoff"
Set-LinkCell $wsZh "C4" "$file2513hash.zh-cn.xlf" "$zhXlfBase/$file2513hash.zh-cn.xlf"
Set-DateCell $wsZh "D4" $newHandoffZh
Set-PlainCell $wsZh "G4" $epoch
Set-PlainCell $wsZh "H4" "Include"

# Row 5 — new file 9834a0f2-...
Set-LinkCell $wsZh "A5" $file9834 "$mdBase/$file9834"
Set-PlainCell $wsZh "B5" "Ready for handoff"
Set-LinkCell $wsZh "C5" "$file9834hash.zh-cn.xlf" "$zhXlfBase/$file9834hash.zh-cn.xlf"
Set-DateCell $wsZh "D5" $newHandoffZh
Set-PlainCell $wsZh "G5" $epoch
Set-PlainCell $wsZh "H5" "Include"

# Row 6 — .localization-config trailing row
Set-LinkCell $wsZh "A6" $cfgFile $cfgUrl
Set-PlainCell $wsZh "B6" "Not to be localized"
Set-DateCell $wsZh "D6" $epoch
Set-PlainCell $wsZh "G6" $epoch
Set-PlainCell $wsZh "H6" "Ignored"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-PlainCell $wsDe "B2" "In Translation"
Set-DateCell  $wsDe "D2" $newHandoffDe

Set-PlainCell $wsDe "B3" "In Translation"
Set-DateCell  $wsDe "D3" $newHandoffDe

# Row 4 — new file 2513cad1-...
Set-LinkCell $wsDe "A4" $file2513 "$mdBase/$file2513"
Set-PlainCell $wsDe "B4" "Ready for handoff"
Set-LinkCell $wsDe "C4" "$file2513hash.de-de.xlf" "$deXlfBase/$file2513hash.de-de.xlf"
Set-DateCell $wsDe "D4" $newHandoffDe
Set-PlainCell $wsDe "G4" $epoch
Set-PlainCell $wsDe "H4" "Include"

# Row 5 — new file 9834a0f2-...
Set-LinkCell $wsDe "A5" $file9834 "$mdBase/$file9834"
Set-PlainCell $wsDe "B5" "Ready for handoff"
Set-LinkCell $wsDe "C5" "$file9834hash.de-de.xlf" "$deXlfBase/$file9834hash.de-de.xlf"
Set-DateCell $wsDe "D5" $newHandoffDe
Set-PlainCell $wsDe "G5" $epoch
Set-PlainCell $wsDe "H5" "Include"

# Row 6 — .localization-config trailing row
Set-LinkCell $wsDe "A6" $cfgFile $cfgUrl
Set-PlainCell $wsDe "B6" "Not to be localized"
Set-DateCell $wsDe "D6" $epoch
Set-PlainCell $wsDe "G6" $epoch
Set-PlainCell $wsDe "H6" "Ignored"

Write-Host "Report regenerated for handoff"
